$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 965.56335
$ws.Range("I15").Value = 965.56335
$ws.Range("K15").Value = 2896.69005
$ws.Range("M15").Value = -2727.69005

$ws.Range("H61").Value = 376.5
$ws.Range("I61").Value = 373.14285
$ws.Range("J61").Value = 400
$ws.Range("K61").Value = 1119.42855
$ws.Range("L61").Value = 1200
$ws.Range("M61").Value = -947.4285500000001
$ws.Range("N61").Value = -1544

$ws.Range("H76").Value = 5584.7144
$ws.Range("I76").Value = 5218.6
$ws.Range("K76").Value = 5218.6
$ws.Range("M76").Value = -4903.6

$ws.Range("H79").Value = 5584.7144
$ws.Range("I79").Value = 5218.6
$ws.Range("K79").Value = 5218.6
$ws.Range("M79").Value = -4126.6

$ws.Range("H118").Value = 998.8823
$ws.Range("I118").Value = 936.3125
$ws.Range("K118").Value = 2808.9375
$ws.Range("M118").Value = -1151.9375

$ws.Range("H127").Value = 1044.6
$ws.Range("I127").Value = 621.25
$ws.Range("J127").Value = 1528.4286
$ws.Range("K127").Value = 1863.75
$ws.Range("L127").Value = 4585.2858
$ws.Range("M127").Value = 3096.25
$ws.Range("N127").Value = -14505.2858

$ws.Range("H131").Value = 765.13336
$ws.Range("I131").Value = 765.13336
$ws.Range("K131").Value = 2295.40008
$ws.Range("M131").Value = 2744.59992

$ws.Range("H132").Value = 1814.9667
$ws.Range("I132").Value = 1401.3928
$ws.Range("K132").Value = 4204.178400000001
$ws.Range("M132").Value = -1674.178400000001

$ws.Range("H137").Value = 377760.53
$ws.Range("I137").Value = 1882.75
$ws.Range("J137").Value = 795402.5
$ws.Range("K137").Value = 5648.25
$ws.Range("L137").Value = 2386207.5
$ws.Range("M137").Value = -3098.25
$ws.Range("N137").Value = -2391307.5

$ws.Range("H138").Value = 2680.0483
$ws.Range("I138").Value = 1897.375
$ws.Range("J138").Value = 3174.3684
$ws.Range("K138").Value = 5692.125
$ws.Range("L138").Value = 9523.1052
$ws.Range("M138").Value = -552.125
$ws.Range("N138").Value = -19803.1052

$ws.Range("H141").Value = 5095.8
$ws.Range("I141").Value = 4550.8887
$ws.Range("K141").Value = 13652.6661
$ws.Range("M141").Value = -8472.666100000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()

$ws.Range("H32").Value = 7002.88
$ws.Range("I32").Value = 4636.841
$ws.Range("J32").Value = 24353.834
$ws.Range("K32").Value = 4636.841
$ws.Range("L32").Value = 24353.834
$ws.Range("M32").Value = -4349.841
$ws.Range("N32").Value = -24927.834

$ws.Range("H74").Value = 40339.777
$ws.Range("I74").Value = 49531.81
$ws.Range("K74").Value = 49531.81
$ws.Range("M74").Value = -48657.81

$ws.Range("H77").Value = 40339.777
$ws.Range("I77").Value = 49531.81
$ws.Range("K77").Value = 247659.05
$ws.Range("M77").Value = -243291.05

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 127102.79
$ws.Range("I20").Value = 183074
$ws.Range("K20").Value = 183074
$ws.Range("M20").Value = -182827

$ws.Range("H134").Value = 2295.724
$ws.Range("I134").Value = 2003.08
$ws.Range("J134").Value = 4124.75
$ws.Range("K134").Value = 6009.24
$ws.Range("L134").Value = 12374.25
$ws.Range("M134").Value = -3474.24
$ws.Range("N134").Value = -17444.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 7539.75
$ws.Range("I4").Value = 53
$ws.Range("J4").Value = 30000
$ws.Range("K4").Value = 53
$ws.Range("L4").Value = 30000
$ws.Range("M4").Value = 59
$ws.Range("N4").Value = -30224

$ws.Range("H7").Value = 14584.857
$ws.Range("I7").Value = 11313.889
$ws.Range("J7").Value = 20472.6
$ws.Range("K7").Value = 11313.889
$ws.Range("L7").Value = 20472.6
$ws.Range("M7").Value = -11200.889
$ws.Range("N7").Value = -20698.6

$ws.Range("H16").Value = 1872.55
$ws.Range("I16").Value = 1763.4
$ws.Range("K16").Value = 1763.4
$ws.Range("M16").Value = -1476.4

$ws.Range("H22").Value = 483.08334
$ws.Range("I22").Value = 483.08334
$ws.Range("K22").Value = 483.08334
$ws.Range("M22").Value = -133.08334

$ws.Range("H31").Value = 4284.7334
$ws.Range("I31").Value = 2268.6155
$ws.Range("J31").Value = 5826.4707
$ws.Range("K31").Value = 2268.6155
$ws.Range("L31").Value = 5826.4707
$ws.Range("M31").Value = -1973.6155
$ws.Range("N31").Value = -6416.4707

$ws.Range("H34").Value = 4284.7334
$ws.Range("I34").Value = 2268.6155
$ws.Range("J34").Value = 5826.4707
$ws.Range("K34").Value = 2268.6155
$ws.Range("L34").Value = 5826.4707
$ws.Range("M34").Value = -2066.6155
$ws.Range("N34").Value = -6230.4707

$ws.Range("H113").Value = 1872.55
$ws.Range("I113").Value = 1763.4
$ws.Range("K113").Value = 1763.4
$ws.Range("M113").Value = 406.5999999999999

$ws.Range("H141").Value = 154320.5
$ws.Range("J141").Value = 165038.58
$ws.Range("L141").Value = 165038.58
$ws.Range("N141").Value = -175398.58

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4032557
$ws.Range("I4").Value = 4310621
$ws.Range("K4").Value = 12931863
$ws.Range("M4").Value = -12931751

$ws.Range("H19").Value = 37.5
$ws.Range("I19").Value = 40
$ws.Range("J19").Value = 30
$ws.Range("K19").Value = 120
$ws.Range("L19").Value = 90
$ws.Range("M19").Value = 54
$ws.Range("N19").Value = -438

$ws.Range("H34").Value = 2026.6666

$ws.Range("H41").Value = 12558.75
$ws.Range("I41").Value = 78
$ws.Range("J41").Value = 50001
$ws.Range("K41").Value = 234
$ws.Range("L41").Value = 150003
$ws.Range("M41").Value = 104
$ws.Range("N41").Value = -150679

$ws.Range("H42").Value = 4416.6665
$ws.Range("J42").Value = 5000
$ws.Range("L42").Value = 15000
$ws.Range("N42").Value = -16068

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 50
$ws.Range("I4").Value = 50
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 50
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 62
$ws.Range("N4").ClearContents()

$ws.Range("H51").Value = 40000
$ws.Range("J51").Value = 40000
$ws.Range("L51").Value = 40000
$ws.Range("N51").Value = -41018

$ws.Range("H113").Value = 2872476.5
$ws.Range("I113").Value = 186518.5
$ws.Range("J113").Value = 5558434.5
$ws.Range("K113").Value = 186518.5
$ws.Range("L113").Value = 5558434.5
$ws.Range("M113").Value = -184348.5
$ws.Range("N113").Value = -5562774.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 9007.154
$ws.Range("I22").Value = 760.75
$ws.Range("J22").Value = 12672.223
$ws.Range("K22").Value = 760.75
$ws.Range("L22").Value = 12672.223
$ws.Range("M22").Value = -465.75
$ws.Range("N22").Value = -13262.223

$ws.Range("H27").Value = 9007.154
$ws.Range("I27").Value = 760.75
$ws.Range("J27").Value = 12672.223
$ws.Range("K27").Value = 760.75
$ws.Range("L27").Value = 12672.223
$ws.Range("M27").Value = -653.75
$ws.Range("N27").Value = -12886.223

$ws.Range("H31").Value = 5432
$ws.Range("J31").Value = 6573.222
$ws.Range("L31").Value = 6573.222
$ws.Range("N31").Value = -7069.222

$ws.Range("H40").Value = 4276251.5
$ws.Range("I40").Value = 2760.8
$ws.Range("J40").Value = 18521220
$ws.Range("K40").Value = 2760.8
$ws.Range("L40").Value = 18521220
$ws.Range("M40").Value = -2624.8
$ws.Range("N40").Value = -18521492

$ws.Range("H122").Value = 20008914
$ws.Range("I122").Value = 9599.714
$ws.Range("K122").Value = 28799.142
$ws.Range("M122").Value = -26349.142

$ws.Range("H132").Value = 2978.077
$ws.Range("I132").Value = 2357.0715
$ws.Range("J132").Value = 3702.5833
$ws.Range("K132").Value = 7071.2145
$ws.Range("L132").Value = 11107.7499
$ws.Range("M132").Value = -4541.2145
$ws.Range("N132").Value = -16167.7499

$ws.Range("H136").Value = 8211.736999999999
$ws.Range("I136").Value = 9095.272000000001
$ws.Range("J136").Value = 6996.875
$ws.Range("K136").Value = 27285.816
$ws.Range("L136").Value = 20990.625
$ws.Range("M136").Value = -24735.816
$ws.Range("N136").Value = -26090.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2894.8572
$ws.Range("I2").Value = 352.8
$ws.Range("K2").Value = 352.8
$ws.Range("M2").Value = -240.8

$ws.Range("H46").Value = 239999
$ws.Range("J46").Value = 239999
$ws.Range("L46").Value = 239999
$ws.Range("N46").Value = -240461

$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()

$ws.Range("H96").Value = 15380.25
$ws.Range("I96").Value = 21460.6
$ws.Range("J96").Value = 5246.3335
$ws.Range("K96").Value = 21460.6
$ws.Range("L96").Value = 5246.3335
$ws.Range("M96").Value = -20087.6
$ws.Range("N96").Value = -7992.3335

$ws.Range("H126").Value = 42579.555
$ws.Range("I126").Value = 47485.125
$ws.Range("J126").Value = 3335
$ws.Range("K126").Value = 142455.375
$ws.Range("L126").Value = 10005
$ws.Range("M126").Value = -139985.375
$ws.Range("N126").Value = -14945

$ws.Range("H134").Value = 239999
$ws.Range("J134").Value = 239999
$ws.Range("L134").Value = 719997
$ws.Range("N134").Value = -725067
